$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.319.05"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.665.47"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.85%  "
$ws.Range("D5").Value = "'219.20"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").Value = "'0.5351"
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("D8").Value = "'0.2665"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("D9").Value = "'0.06402"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D10").Value = "'20.72"
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("D11").Value = "'0.07843"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "'4.571"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").Value = "1.661.51"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").Value = "1.893.84"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "'0.5531"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "0.0₅8220"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "'1.011"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "'4.693"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").Value = "'193.79"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").Value = "'6.044"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").Value = "'146.40"
$ws.Range("E24").Value = "  +2.91%  "
$ws.Range("D25").Value = "'0.1233"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "'7.209"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "'1.502"
$ws.Range("E28").Value = "  +5.09%  "
$ws.Range("D29").Value = "'0.05841"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").Value = "'1.282"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").Value = "'3.653"
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").Value = "'1.616"
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("D34").Value = "'0.9699"
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("D35").Value = "'2.825"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("D36").Value = "'2.421"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").Value = "'0.5824"
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("D38").Value = "'0.01607"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "'0.8748"
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("D40").Value = "'5.866"
$ws.Range("E40").Value = "  +2.11%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.052.95"
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'105.28"
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "1.804.79"
$ws.Range("D45").Value = "'57.92"
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("D46").Value = "'1.014"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈104"
$ws.Range("E47").Value = "  -7.38%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.4388"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").Value = "'8.007"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").Value = "'0.05168"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").Value = "'1.414"
$ws.Range("E51").Value = "  -3.53%  "
